$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.228521823883057
$ws.Range("B1").Value = 1.722233533859253
$ws.Range("C1").Value = 3.101752519607544
$ws.Range("D1").Value = 3.745893955230713
$ws.Range("E1").Value = 1.369771838188171
